# ZBP_06_home_office.xlsx update: add new weekly survey columns
# (week of 4.–10. 10. 2021 and 11.–17. 10. 2021) to both the "data" and
# "pocetR" sheets, and bump the "aktualizace" date stamp in each sheet's
# final footer row from 6. 10. 2021 to 20. 10. 2021.
#
# NOTE: the two new week-label strings below are embedded as literal
# text (en dash U+2013 typed directly into this UTF-8 file) rather than
# built with string concatenation / [char] casts, since this host's
# PowerShell-emulation "+" operator coerces a numeric-looking string
# plus a char into numeric addition instead of concatenating text.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "data": new header date labels in BN1 / BO1, new data columns
# BN (rows 2-77) and BO (rows 2-77).
# ---------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("data")

$bnHeader = "4.–10. 10. 2021"
$boHeader = "11.–17. 10. 2021"

$headerRng = $wsData.Range("BN1:BO1")
$headerArr = New-Object 'object[,]' 1,2
$headerArr[0,0] = $bnHeader
$headerArr[0,1] = $boHeader
$headerRng.Value = $headerArr

# Match the look of the existing header cells (bold, centered, thin box
# border) used by the rest of row 1.
$headerRng.Font.Bold = $true
$headerRng.HorizontalAlignment = -4108
$headerRng.VerticalAlignment = -4160
$headerRng.Borders.LineStyle = 1

$bnVals = @(0.78,0.08,0.06,0.08,0.78,0.07000000000000001,0.07000000000000001,0.08,0.45,0.25,0.27,0.03,0.84,0.05,0.02,0.09,0.93,0.01,0.005,0.055,0.8,0.14,0.005,0.055,0.82,0.07000000000000001,0.05,0.06,0.8100000000000001,0.05,0.04,0.1,0.73,0.11,0.1,0.06,0.64,0.15,0.13,0.08,0.78,0.09,0.09,0.04,0.84,0.07000000000000001,0.05,0.04,0.61,0.11,0.09,0.19,0.77,0.05,0.04,0.14,0.74,0.15,0.07000000000000001,0.04,0.8,0.07000000000000001,0.07000000000000001,0.06,0.84,0.04,0.02,0.1,0.78,0.09,0.08,0.05,0.7,0.13,0.11,0.06)
$boVals = @(0.78,0.07000000000000001,0.07000000000000001,0.08,0.75,0.07000000000000001,0.08,0.1,0.47,0.18,0.31,0.04,0.84,0.04,0.02,0.1,0.89,0.01,0.005,0.095,0.82,0.1,0.005,0.075,0.8100000000000001,0.05,0.06,0.08,0.8,0.05,0.04,0.11,0.73,0.09,0.1,0.08,0.66,0.13,0.13,0.08,0.73,0.11,0.08,0.08,0.84,0.05,0.06,0.05,0.64,0.08,0.1,0.18,0.76,0.04,0.06,0.14,0.74,0.13,0.1,0.03,0.79,0.06,0.07000000000000001,0.08,0.82,0.03,0.03,0.12,0.78,0.09,0.08,0.05,0.71,0.1,0.13,0.06)

$dataRng = $wsData.Range("BN2:BO77")
$dataArr = New-Object 'object[,]' 76,2
for ($i = 0; $i -lt 76; $i++) {
    $dataArr[$i,0] = $bnVals[$i]
    $dataArr[$i,1] = $boVals[$i]
}
$dataRng.Value = $dataArr

# Footer row: bump the "aktualizace" date stamp.
$wsData.Range("A78").Value = "Život během pandemie, Home office, % respondentů celkově a ve skupinách, aktualizace 20. 10. 2021"

# ---------------------------------------------------------------------
# Sheet "pocetR": new header date labels in BM1 / BN1, new data columns
# BM (rows 2-20) and BN (rows 2-20).
# ---------------------------------------------------------------------
$wsPocet = $wb.Worksheets.Item("pocetR")

$headerRng2 = $wsPocet.Range("BM1:BN1")
$headerArr2 = New-Object 'object[,]' 1,2
$headerArr2[0,0] = $bnHeader
$headerArr2[0,1] = $boHeader
$headerRng2.Value = $headerArr2

$headerRng2.Font.Bold = $true
$headerRng2.HorizontalAlignment = -4108
$headerRng2.VerticalAlignment = -4160
$headerRng2.Borders.LineStyle = 1

$bmVals2 = @(1041,304,101,258,150,85,518,245,131,147,265,589,187,191,164,686,390,411,240)
$bnVals2 = @(1041,304,101,258,150,85,518,245,131,147,265,589,187,191,164,686,390,411,240)

$dataRng2 = $wsPocet.Range("BM2:BN20")
$dataArr2 = New-Object 'object[,]' 19,2
for ($i = 0; $i -lt 19; $i++) {
    $dataArr2[$i,0] = $bmVals2[$i]
    $dataArr2[$i,1] = $bnVals2[$i]
}
$dataRng2.Value = $dataArr2

# Footer row: bump the "aktualizace" date stamp.
$wsPocet.Range("A21").Value = "Život během pandemie, Home office, velikost dotázaného souboru celkově a ve skupinách, aktualizace 20. 10. 2021"
